$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.038726317646665
$ws.Range("D2").Value2 = 1.041321537703689
$ws.Range("E2").Value2 = 1.056370576205585
$ws.Range("F2").Value2 = 1.063060475827773
$ws.Range("I2").Value2 = 1.042593218532114
$ws.Range("J2").Value2 = 1.043822191797132
$ws.Range("K2").Value2 = 1.044101382770303
$ws.Range("L2").Value2 = 1.059108472584678
$ws.Range("M2").Value2 = 1.065780137225344
$ws.Range("N2").Value2 = 1.018510623207622
$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.03968169000082
$ws.Range("D3").Value2 = 1.042038172224527
$ws.Range("E3").Value2 = 1.057484924250156
$ws.Range("F3").Value2 = 1.06420815509101
$ws.Range("I3").Value2 = 1.042884144298797
$ws.Range("J3").Value2 = 1.044422446679903
$ws.Range("K3").Value2 = 1.044628868010215
$ws.Range("L3").Value2 = 1.06003568074297
$ws.Range("M3").Value2 = 1.066741918185743
$ws.Range("N3").Value2 = 1.018712133983734
$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.04029984720686
$ws.Range("D4").Value2 = 1.042501809280626
$ws.Range("E4").Value2 = 1.058206903004421
$ws.Range("F4").Value2 = 1.064951550039955
$ws.Range("I4").Value2 = 1.043071061864873
$ws.Range("J4").Value2 = 1.044810202944001
$ws.Range("K4").Value2 = 1.044969442007843
$ws.Range("L4").Value2 = 1.060635976154611
$ws.Range("M4").Value2 = 1.067364438569524
$ws.Range("N4").Value2 = 1.018842246554351
$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.04055971230747
$ws.Range("D5").Value2 = 1.042696703845148
$ws.Range("E5").Value2 = 1.058510643291825
$ws.Range("F5").Value2 = 1.065264257458632
$ws.Range("I5").Value2 = 1.043149323078024
$ws.Range("J5").Value2 = 1.044973059872502
$ws.Range("K5").Value2 = 1.045112440605457
$ws.Range("L5").Value2 = 1.060888419193333
$ws.Range("M5").Value2 = 1.067626189999342
$ws.Range("N5").Value2 = 1.018896879008377
$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.040603344309946
$ws.Range("D6").Value2 = 1.042729426392286
$ws.Range("E6").Value2 = 1.058561655545792
$ws.Range("F6").Value2 = 1.065316773208215
$ws.Range("I6").Value2 = 1.04316244476845
$ws.Range("J6").Value2 = 1.045000395116837
$ws.Range("K6").Value2 = 1.045136440203758
$ws.Range("L6").Value2 = 1.060930810111469
$ws.Range("M6").Value2 = 1.067670141784228
$ws.Range("N6").Value2 = 1.018906048109878
$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.040303319569706
$ws.Range("D7").Value2 = 1.042504413544958
$ws.Range("E7").Value2 = 1.058210960731003
$ws.Range("F7").Value2 = 1.064955727728181
$ws.Range("I7").Value2 = 1.043072108848617
$ws.Range("J7").Value2 = 1.044812379659064
$ws.Range("K7").Value2 = 1.044971353464996
$ws.Range("L7").Value2 = 1.060639349001138
$ws.Range("M7").Value2 = 1.067367935932764
$ws.Range("N7").Value2 = 1.018842976818916
$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.039049196722472
$ws.Range("D8").Value2 = 1.041563742363562
$ws.Range("E8").Value2 = 1.056746984582709
$ws.Range("F8").Value2 = 1.063448179976226
$ws.Range("I8").Value2 = 1.042691813626883
$ws.Range("J8").Value2 = 1.044025184788609
$ws.Range("K8").Value2 = 1.044279802647854
$ws.Range("L8").Value2 = 1.059421758593843
$ws.Range("M8").Value2 = 1.066105137464375
$ws.Range("N8").Value2 = 1.018578782221269
$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.036839039112034
$ws.Range("D9").Value2 = 1.039905636051134
$ws.Range("E9").Value2 = 1.054174329147056
$ws.Range("F9").Value2 = 1.060797587073086
$ws.Range("I9").Value2 = 1.042011506098652
$ws.Range("J9").Value2 = 1.0426330990818
$ws.Range("K9").Value2 = 1.043055525559526
$ws.Range("L9").Value2 = 1.057278737856202
$ws.Range("M9").Value2 = 1.063881338690812
$ws.Range("N9").Value2 = 1.018111115478463
$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.035365464114954
$ws.Range("D10").Value2 = 1.038799930512709
$ws.Range("E10").Value2 = 1.052463972546998
$ws.Range("F10").Value2 = 1.059034493304574
$ws.Range("I10").Value2 = 1.041551141986512
$ws.Range("J10").Value2 = 1.041701742215055
$ws.Range("K10").Value2 = 1.042235555002988
$ws.Range("L10").Value2 = 1.055851760521005
$ws.Range("M10").Value2 = 1.062399761902335
$ws.Range("N10").Value2 = 1.017797923729422
$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.034727360880945
$ws.Range("D11").Value2 = 1.038321086547767
$ws.Range("E11").Value2 = 1.051724494709038
$ws.Range("F11").Value2 = 1.058271994940317
$ws.Range("I11").Value2 = 1.041350183672612
$ws.Range("J11").Value2 = 1.041297676889473
$ws.Range("K11").Value2 = 1.041879606723164
$ws.Range("L11").Value2 = 1.055234267338833
$ws.Range("M11").Value2 = 1.061758449790927
$ws.Range("N11").Value2 = 1.017661975424502
$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.034490335762843
$ws.Range("D12").Value2 = 1.038143213403819
$ws.Range("E12").Value2 = 1.051449987519375
$ws.Range("F12").Value2 = 1.057988909100592
$ws.Range("I12").Value2 = 1.041275295994724
$ws.Range("J12").Value2 = 1.041147471942907
$ws.Range("K12").Value2 = 1.041747257544172
$ws.Range("L12").Value2 = 1.055004962472434
$ws.Range("M12").Value2 = 1.061520270955365
$ws.Range("N12").Value2 = 1.01761142814097
$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.034541178686223
$ws.Range("D13").Value2 = 1.038181368199902
$ws.Range("E13").Value2 = 1.051508862652092
$ws.Range("F13").Value2 = 1.058049625641109
$ws.Range("I13").Value2 = 1.041291370640139
$ws.Range("J13").Value2 = 1.041179696707081
$ws.Range("K13").Value2 = 1.041775652955728
$ws.Range("L13").Value2 = 1.055054146449617
$ws.Range("M13").Value2 = 1.061571359630321
$ws.Range("N13").Value2 = 1.017622272965927
$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.034707768415105
$ws.Range("D14").Value2 = 1.038306383680941
$ws.Range("E14").Value2 = 1.05170180043674
$ws.Range("F14").Value2 = 1.058248592121331
$ws.Range("I14").Value2 = 1.041343998385843
$ws.Range("J14").Value2 = 1.041285263288434
$ws.Range("K14").Value2 = 1.041868669432731
$ws.Range("L14").Value2 = 1.055215311711952
$ws.Range("M14").Value2 = 1.061738761180568
$ws.Range("N14").Value2 = 1.017657798189826
$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.034810409165976
$ws.Range("D15").Value2 = 1.03838340865197
$ws.Range("E15").Value2 = 1.051820698005232
$ws.Range("F15").Value2 = 1.058371200483858
$ws.Range("I15").Value2 = 1.041376391908382
$ws.Range("J15").Value2 = 1.041350290837513
$ws.Range("K15").Value2 = 1.041925962191861
$ws.Range("L15").Value2 = 1.055314618842877
$ws.Range("M15").Value2 = 1.061841907176101
$ws.Range("N15").Value2 = 1.017679679828713
$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.035407812191113
$ws.Range("D16").Value2 = 1.038831708480981
$ws.Range("E16").Value2 = 1.052513072795522
$ws.Range("F16").Value2 = 1.059085117451373
$ws.Range("I16").Value2 = 1.04156444486508
$ws.Range("J16").Value2 = 1.041728542260495
$ws.Range("K16").Value2 = 1.042259159284016
$ws.Range("L16").Value2 = 1.05589274987332
$ws.Range("M16").Value2 = 1.062442328312208
$ws.Range("N16").Value2 = 1.017806939144737
$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.035782537933672
$ws.Range("D17").Value2 = 1.039112897946149
$ws.Range("E17").Value2 = 1.052947680390455
$ws.Range("F17").Value2 = 1.059533188352548
$ws.Range("I17").Value2 = 1.041681972587189
$ws.Range("J17").Value2 = 1.041965600354395
$ws.Range("K17").Value2 = 1.04246792552944
$ws.Range("L17").Value2 = 1.05625550251475
$ws.Range("M17").Value2 = 1.062819015583259
$ws.Range("N17").Value2 = 1.017886676201999
$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.036001105679033
$ws.Range("D18").Value2 = 1.039276904568118
$ws.Range("E18").Value2 = 1.053201287513733
$ws.Range("F18").Value2 = 1.059794630546939
$ws.Range("I18").Value2 = 1.04175036842812
$ws.Range("J18").Value2 = 1.042103796888488
$ws.Range("K18").Value2 = 1.042589608906279
$ws.Range("L18").Value2 = 1.056467128345846
$ws.Range("M18").Value2 = 1.063038752167263
$ws.Range("N18").Value2 = 1.017933153219812
$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.036075631043089
$ws.Range("D19").Value2 = 1.039332825472928
$ws.Range("E19").Value2 = 1.053287779289289
$ws.Range("F19").Value2 = 1.059883790852294
$ws.Range("I19").Value2 = 1.041773663160077
$ws.Range("J19").Value2 = 1.04215090552874
$ws.Range("K19").Value2 = 1.042631085118422
$ws.Range("L19").Value2 = 1.056539293794111
$ws.Range("M19").Value2 = 1.063113680274407
$ws.Range("N19").Value2 = 1.017948995204663
$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.035742333762371
$ws.Range("D20").Value2 = 1.039082729637766
$ws.Range("E20").Value2 = 1.052901039952411
$ws.Range("F20").Value2 = 1.059485105275286
$ws.Range("I20").Value2 = 1.041669379118496
$ws.Range("J20").Value2 = 1.041940174065138
$ws.Range("K20").Value2 = 1.042445535825087
$ws.Range("L20").Value2 = 1.056216578623103
$ws.Range("M20").Value2 = 1.062778598404219
$ws.Range("N20").Value2 = 1.01787812450571
$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.034658712033279
$ws.Range("D21").Value2 = 1.038269569984831
$ws.Range("E21").Value2 = 1.051644980436885
$ws.Range("F21").Value2 = 1.058189997588431
$ws.Range("I21").Value2 = 1.041328507527517
$ws.Range("J21").Value2 = 1.041254179791018
$ws.Range("K21").Value2 = 1.041841282098181
$ws.Range("L21").Value2 = 1.05516785091489
$ws.Range("M21").Value2 = 1.061689464685013
$ws.Range("N21").Value2 = 1.017647338276447
$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.033977366244956
$ws.Range("D22").Value2 = 1.037758251531657
$ws.Range("E22").Value2 = 1.050856216766317
$ws.Range("F22").Value2 = 1.057376521926207
$ws.Range("I22").Value2 = 1.041112783298956
$ws.Range("J22").Value2 = 1.040822190665279
$ws.Range("K22").Value2 = 1.041460587915071
$ws.Range("L22").Value2 = 1.054508818928777
$ws.Range("M22").Value2 = 1.061004874194899
$ws.Range("N22").Value2 = 1.017501944445765
$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.034338563305031
$ws.Range("D23").Value2 = 1.038029315858732
$ws.Range("E23").Value2 = 1.051274263337134
$ws.Range("F23").Value2 = 1.05780768404204
$ws.Range("I23").Value2 = 1.041227275873965
$ws.Range("J23").Value2 = 1.041051260360432
$ws.Range("K23").Value2 = 1.041662474510036
$ws.Range("L23").Value2 = 1.054858151612245
$ws.Range("M23").Value2 = 1.061367770476799
$ws.Range("N23").Value2 = 1.017579047834646
$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.035760500307346
$ws.Range("D24").Value2 = 1.039096361417658
$ws.Range("E24").Value2 = 1.052922114425306
$ws.Range("F24").Value2 = 1.05950683166988
$ws.Range("I24").Value2 = 1.041675070047183
$ws.Range("J24").Value2 = 1.041951663343853
$ws.Range("K24").Value2 = 1.042455653035802
$ws.Range("L24").Value2 = 1.056234166535636
$ws.Range("M24").Value2 = 1.06279686112182
$ws.Range("N24").Value2 = 1.01788198874883
$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.037410443374558
$ws.Range("D25").Value2 = 1.040334353085989
$ws.Range("E25").Value2 = 1.054838585703012
$ws.Range("F25").Value2 = 1.061482129625642
$ws.Range("I25").Value2 = 1.042188586298147
$ws.Range("J25").Value2 = 1.042993570170082
$ws.Range("K25").Value2 = 1.043372700041672
$ws.Range("L25").Value2 = 1.057832460104854
$ws.Range("M25").Value2 = 1.064456075931851
$ws.Range("N25").Value2 = 1.01823226869471
